$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 38255.25
$ws.Range("I20").Value = 38255.25
$ws.Range("K20").Value = 38255.25
$ws.Range("M20").Value = -38025.25
$ws.Range("H33").Value = 111.23077
$ws.Range("I33").Value = 103.875
$ws.Range("K33").Value = 103.875
$ws.Range("M33").Value = 125.125
$ws.Range("H35").Value = 38255.25
$ws.Range("I35").Value = 38255.25
$ws.Range("K35").Value = 38255.25
$ws.Range("M35").Value = -37876.25
$ws.Range("H64").Value = 4507.85
$ws.Range("I64").Value = 3300.1538
$ws.Range("J64").Value = 6750.7144
$ws.Range("K64").Value = 3300.1538
$ws.Range("L64").Value = 6750.7144
$ws.Range("M64").Value = -3052.1538
$ws.Range("N64").Value = -7246.7144
$ws.Range("H67").Value = 4507.85
$ws.Range("I67").Value = 3300.1538
$ws.Range("J67").Value = 6750.7144
$ws.Range("K67").Value = 3300.1538
$ws.Range("L67").Value = 6750.7144
$ws.Range("M67").Value = -2442.1538
$ws.Range("N67").Value = -8466.714400000001
$ws.Range("H70").Value = 2520.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2520.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 7561.5
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -8101.5
$ws.Range("H73").Value = 2520.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2520.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 7561.5
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -9433.5
$ws.Range("H76").Value = 2587433.5
$ws.Range("I76").Value = 3370365
$ws.Range("J76").Value = 3759.8
$ws.Range("K76").Value = 3370365
$ws.Range("L76").Value = 3759.8
$ws.Range("M76").Value = -3370050
$ws.Range("N76").Value = -4389.8
$ws.Range("H79").Value = 2587433.5
$ws.Range("I79").Value = 3370365
$ws.Range("J79").Value = 3759.8
$ws.Range("K79").Value = 3370365
$ws.Range("L79").Value = 3759.8
$ws.Range("M79").Value = -3369273
$ws.Range("N79").Value = -5943.8
$ws.Range("H132").Value = 557096.8
$ws.Range("I132").Value = 675134.7
$ws.Range("K132").Value = 2025404.1
$ws.Range("M132").Value = -2022874.1
$ws.Range("H133").Value = 43958.184
$ws.Range("J133").Value = 43958.184
$ws.Range("L133").Value = 43958.184
$ws.Range("N133").Value = -54078.184
$ws.Range("H138").Value = 6929290.5
$ws.Range("J138").Value = 8335964
$ws.Range("L138").Value = 25007892
$ws.Range("N138").Value = -25018172
$ws.Range("H139").Value = 44750
$ws.Range("J139").Value = 44750
$ws.Range("L139").Value = 44750
$ws.Range("N139").Value = -55030
$ws.Range("H141").Value = 3910.5557
$ws.Range("I141").Value = 2653
$ws.Range("K141").Value = 7959
$ws.Range("M141").Value = -2779

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28098.215
$ws.Range("I32").Value = 3877.775
$ws.Range("K32").Value = 3877.775
$ws.Range("M32").Value = -3590.775
$ws.Range("H54").Value = 6300
$ws.Range("I54").Value = 3000
$ws.Range("K54").Value = 3000
$ws.Range("M54").Value = -2231
$ws.Range("H61").Value = 2836.75
$ws.Range("I61").Value = 2260.24
$ws.Range("J61").Value = 4895.7144
$ws.Range("K61").Value = 2260.24
$ws.Range("L61").Value = 4895.7144
$ws.Range("M61").Value = -2048.24
$ws.Range("N61").Value = -5319.7144
$ws.Range("H74").Value = 4966.5
$ws.Range("I74").Value = 886.2593000000001
$ws.Range("K74").Value = 886.2593000000001
$ws.Range("M74").Value = -12.25930000000005
$ws.Range("H77").Value = 4966.5
$ws.Range("I77").Value = 886.2593000000001
$ws.Range("K77").Value = 4431.2965
$ws.Range("M77").Value = -63.29650000000038
$ws.Range("H97").Value = 30303590
$ws.Range("I97").Value = 37037420
$ws.Range("J97").Value = 1355.5
$ws.Range("K97").Value = 37037420
$ws.Range("L97").Value = 1355.5
$ws.Range("M97").Value = -37036924
$ws.Range("N97").Value = -2347.5
$ws.Range("H125").Value = 31737.625
$ws.Range("I125").Value = 30000
$ws.Range("J125").Value = 31985.857
$ws.Range("K125").Value = 30000
$ws.Range("L125").Value = 31985.857
$ws.Range("M125").Value = -25080
$ws.Range("N125").Value = -41825.857
$ws.Range("H132").Value = 2972.6458
$ws.Range("I132").Value = 2569.718
$ws.Range("J132").Value = 4718.6665
$ws.Range("K132").Value = 7709.154
$ws.Range("L132").Value = 14155.9995
$ws.Range("M132").Value = -5179.154
$ws.Range("N132").Value = -19215.9995
$ws.Range("H136").Value = 2836.75
$ws.Range("I136").Value = 2260.24
$ws.Range("J136").Value = 4895.7144
$ws.Range("K136").Value = 6780.719999999999
$ws.Range("L136").Value = 14687.1432
$ws.Range("M136").Value = -4230.719999999999
$ws.Range("N136").Value = -19787.1432
$ws.Range("H139").Value = 37195
$ws.Range("J139").Value = 37195
$ws.Range("L139").Value = 37195
$ws.Range("N139").Value = -47475

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H94").Value = 1191.6522
$ws.Range("I94").Value = 1028.9524
$ws.Range("J94").Value = 2900
$ws.Range("K94").Value = 1028.9524
$ws.Range("L94").Value = 2900
$ws.Range("M94").Value = -577.9523999999999
$ws.Range("N94").Value = -3802
$ws.Range("H96").Value = 30000
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -35492
$ws.Range("H134").Value = 3958.8076
$ws.Range("I134").Value = 2946.8667
$ws.Range("J134").Value = 5338.727
$ws.Range("K134").Value = 8840.6001
$ws.Range("L134").Value = 16016.181
$ws.Range("M134").Value = -6305.6001
$ws.Range("N134").Value = -21086.181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2165.4
$ws.Range("J31").Value = 3428.1875
$ws.Range("L31").Value = 3428.1875
$ws.Range("N31").Value = -4018.1875
$ws.Range("H34").Value = 2165.4
$ws.Range("J34").Value = 3428.1875
$ws.Range("L34").Value = 3428.1875
$ws.Range("N34").Value = -3832.1875
$ws.Range("H93").Value = 10409.777
$ws.Range("I93").Value = 7961
$ws.Range("J93").Value = 30000
$ws.Range("K93").Value = 7961
$ws.Range("L93").Value = 30000
$ws.Range("M93").Value = -6089
$ws.Range("N93").Value = -33744
$ws.Range("H134").Value = 5101.6665
$ws.Range("I134").Value = 2704.8
$ws.Range("K134").Value = 8114.400000000001
$ws.Range("M134").Value = -5579.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 972
$ws.Range("I5").Value = 561.5294
$ws.Range("J5").Value = 1669.8
$ws.Range("K5").Value = 1684.5882
$ws.Range("L5").Value = 5009.4
$ws.Range("M5").Value = -1572.5882
$ws.Range("N5").Value = -5233.4
$ws.Range("H58").Value = 8251.143
$ws.Range("J58").Value = 9166.666999999999
$ws.Range("L58").Value = 27500.001
$ws.Range("N58").Value = -27756.001
$ws.Range("H135").Value = 972
$ws.Range("I135").Value = 561.5294
$ws.Range("J135").Value = 1669.8
$ws.Range("K135").Value = 5053.7646
$ws.Range("L135").Value = 15028.2
$ws.Range("M135").Value = -2518.7646
$ws.Range("N135").Value = -20098.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H97").Value = 1255.25
$ws.Range("I97").Value = 1003.3333
$ws.Range("J97").Value = 2011
$ws.Range("K97").Value = 1003.3333
$ws.Range("L97").Value = 2011
$ws.Range("M97").Value = -507.3333
$ws.Range("N97").Value = -3003
$ws.Range("H132").Value = 3479.625
$ws.Range("I132").Value = 3496.923
$ws.Range("K132").Value = 10490.769
$ws.Range("M132").Value = -7960.769
$ws.Range("H138").Value = 64532.08
$ws.Range("J138").Value = 64532.08
$ws.Range("L138").Value = 64532.08
$ws.Range("N138").Value = -74812.08

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1932.4546
$ws.Range("I68").Value = 1824
$ws.Range("J68").Value = 2122.25
$ws.Range("K68").Value = 1824
$ws.Range("L68").Value = 2122.25
$ws.Range("M68").Value = -1075
$ws.Range("N68").Value = -3620.25
$ws.Range("H71").Value = 1932.4546
$ws.Range("I71").Value = 1824
$ws.Range("J71").Value = 2122.25
$ws.Range("K71").Value = 9120
$ws.Range("L71").Value = 10611.25
$ws.Range("M71").Value = -5376
$ws.Range("N71").Value = -18099.25
$ws.Range("H93").Value = 4056.7144
$ws.Range("I93").Value = 3680.6
$ws.Range("J93").Value = 4997
$ws.Range("K93").Value = 3680.6
$ws.Range("L93").Value = 4997
$ws.Range("M93").Value = -2432.6
$ws.Range("N93").Value = -7493

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4750
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 4750
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 4750
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -7496
$ws.Range("H136").Value = 1425.58
$ws.Range("I136").Value = 923.5925999999999
$ws.Range("K136").Value = 2770.7778
$ws.Range("M136").Value = -220.7777999999998
